# Repull data, push all data, mean calculation
# Update the dSF column (F) with freshly pulled values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 5
    3  = -2
    4  = -6
    6  = 3
    7  = 1
    8  = 0
    9  = 4
    10 = 1
    11 = 2
    12 = 8
    13 = 6
    14 = 3
    15 = -6
    16 = -2
    17 = 9
    18 = 6
    20 = 2
    21 = -5
    23 = -4
    24 = 3
    25 = 2
    27 = 2
    28 = 4
    29 = 1
    30 = 2
    32 = 5
    33 = -1
    34 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
